$wb = $excel.ActiveWorkbook

# "Komplett und Minderest GW" - fill the GW_min (E) and GW_komp (F) columns on the
# GWK sheet with the constant lookup keys G001 / G037 for every weight-bracket row.
$wsGWK = $wb.Worksheets.Item("GWK")
$wsGWK.Range("E2:E36").Value = "G001"
$wsGWK.Range("F2:F36").Value = "G037"

# Update the current selection on GWK (user ended up with D40 selected there) ...
$wsGWK.Activate()
$wsGWK.Range("D40").Select()

# ... then switch back to the Zonen sheet, which becomes the active/visible tab.
$wsZonen = $wb.Worksheets.Item("Zonen")
$wsZonen.Activate()
